$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'basketball gear boys'
$ws.Cells.Item(2, 1).Value = 'spandex shorts for volleyball'
$ws.Cells.Item(3, 1).Value = 'volleyball knee sleeves men'
$ws.Cells.Item(4, 1).Value = 'knee pads for work construction'
$ws.Cells.Item(5, 1).Value = 'knee pads for working on floors'
$ws.Cells.Item(6, 1).Value = 'bee tights adult'
$ws.Cells.Item(7, 1).Value = 'compression test'
$ws.Cells.Item(8, 1).Value = 'lacrosse youth shorts'
$ws.Cells.Item(9, 1).Value = 'women volleyball knee pads'
$ws.Cells.Item(10, 1).Value = 'short football pants'
$ws.Cells.Item(11, 1).Value = 'girls compression tights'
$ws.Cells.Item(12, 1).Value = 'knee compression sleeve reduce strain & swelling'
$ws.Cells.Item(13, 1).Value = 'volleyball knee pads small'
$ws.Cells.Item(14, 1).Value = 'yoga pants knee length'
$ws.Cells.Item(15, 1).Value = 'pad for squats'
$ws.Cells.Item(16, 1).Value = 'sweat pads'
$ws.Cells.Item(17, 1).Value = 'calf protector'
$ws.Cells.Item(18, 1).Value = 'compression pants long'
$ws.Cells.Item(19, 1).Value = 'spandex mens pants'
$ws.Cells.Item(20, 1).Value = 'protective knee pads for work'
$ws.Cells.Item(21, 1).Value = 'youth compression shorts'
$ws.Cells.Item(22, 1).Value = 'mens basketball shorts long'
$ws.Cells.Item(23, 1).Value = 'thick volleyball knee pads'
$ws.Cells.Item(24, 1).Value = 'weightlifting guide'
$ws.Cells.Item(25, 1).Value = 'small work knee pads'
$ws.Cells.Item(26, 1).Value = 'basketballs in bulk'
$ws.Cells.Item(27, 1).Value = 'baseballs cheap'
$ws.Cells.Item(28, 1).Value = 'spandex capri leggings'
$ws.Cells.Item(29, 1).Value = 'hex fabric'
$ws.Cells.Item(30, 1).Value = 'black knee pads for work'
$ws.Cells.Item(31, 1).Value = 'black knee guards'
$ws.Cells.Item(32, 1).Value = 'work knee pads for men gel'
$ws.Cells.Item(33, 1).Value = 'leg guard baseball'
$ws.Cells.Item(34, 1).Value = 'lacrosse shorts youth'
$ws.Cells.Item(35, 1).Value = 'girdles for men'
$ws.Cells.Item(36, 1).Value = 'basketball cheap'
$ws.Cells.Item(37, 1).Value = 'wrestling fight shorts'
$ws.Cells.Item(38, 1).Value = 'construction knee'
$ws.Cells.Item(39, 1).Value = 'knee compression sleeve protector'
$ws.Cells.Item(40, 1).Value = 'knee bursitis sleeve'
$ws.Cells.Item(41, 1).Value = 'football leggings'
$ws.Cells.Item(42, 1).Value = 'black girls softball pants'
$ws.Cells.Item(43, 1).Value = 'girls softball pants black'
$ws.Cells.Item(44, 1).Value = 'working pants with knee pads'
$ws.Cells.Item(45, 1).Value = 'youth volleyball sleeves'
$ws.Cells.Item(46, 1).Value = 'mens knee length shorts'
$ws.Cells.Item(47, 1).Value = 'compression capris girls'
$ws.Cells.Item(48, 1).Value = 'leaf leggings'
$ws.Cells.Item(49, 1).Value = 'fit compression knee'
$ws.Cells.Item(50, 1).Value = 'performance basketball'
$ws.Cells.Item(51, 1).Value = 'pantalones de basketball'
$ws.Cells.Item(52, 1).Value = 'youth baseball sleeves for boys'
$ws.Cells.Item(53, 1).Value = 'mens tights for sports'
$ws.Cells.Item(54, 1).Value = 'knee protector work'
$ws.Cells.Item(55, 1).Value = 'long basketball shorts for men'
$ws.Cells.Item(56, 1).Value = 'mens knee sleeves weightlifting'
$ws.Cells.Item(57, 1).Value = 'knee pads bulk'
$ws.Cells.Item(58, 1).Value = 'snowboarding protective gear'
$ws.Cells.Item(59, 1).Value = 'calf compression sleeve youth'
$ws.Cells.Item(60, 1).Value = 'knee pain pads'
$ws.Cells.Item(61, 1).Value = 'knee pad for construction'
$ws.Cells.Item(62, 1).Value = 'mens capri yoga pants'
$ws.Cells.Item(63, 1).Value = 'knee sleeves for basketball'
$ws.Cells.Item(64, 1).Value = 'knee pad sleeves'
$ws.Cells.Item(65, 1).Value = 'knee work'
$ws.Cells.Item(66, 1).Value = 'volleyball shorts men'
$ws.Cells.Item(67, 1).Value = 'knee pad for working'
$ws.Cells.Item(68, 1).Value = 'above knee shorts men'
$ws.Cells.Item(69, 1).Value = 'youth volleyball shorts for girls'
$ws.Cells.Item(70, 1).Value = 'protector paintball'
$ws.Cells.Item(71, 1).Value = 'baseball shorts boys'
$ws.Cells.Item(72, 1).Value = 'boys tights and leggings'
$ws.Cells.Item(73, 1).Value = 'kneeling pad gel'
$ws.Cells.Item(74, 1).Value = 'paintball pants men'
$ws.Cells.Item(75, 1).Value = 'men above knee shorts'
$ws.Cells.Item(76, 1).Value = 'spandex compression shorts'
$ws.Cells.Item(77, 1).Value = 'hockey padded shorts'
$ws.Cells.Item(78, 1).Value = 'knee pads for joint pain'
$ws.Cells.Item(79, 1).Value = 'rodillera volleyball'
$ws.Cells.Item(80, 1).Value = 'basketball calf sleeve'
$ws.Cells.Item(81, 1).Value = 'womens lacrosse pants'
$ws.Cells.Item(82, 1).Value = 'best work knee pads'
$ws.Cells.Item(83, 1).Value = 'leg sleeves for men basketball'
$ws.Cells.Item(84, 1).Value = 'football pouch youth'
$ws.Cells.Item(85, 1).Value = 'calf compression sleeve boys'
$ws.Cells.Item(86, 1).Value = 'boys leggings sports'
$ws.Cells.Item(87, 1).Value = 'basketball stretch pants'
$ws.Cells.Item(88, 1).Value = 'mens spandex shorts'
$ws.Cells.Item(89, 1).Value = 'black kneepads'
$ws.Cells.Item(90, 1).Value = 'joint protectors'
$ws.Cells.Item(91, 1).Value = 'boys hiking pants'
$ws.Cells.Item(92, 1).Value = 'knee length yoga pants'
$ws.Cells.Item(93, 1).Value = 'girls volleyball shorts youth'
$ws.Cells.Item(94, 1).Value = 'baseball pants mens long'
$ws.Cells.Item(95, 1).Value = 'mens sports tights'
$ws.Cells.Item(96, 1).Value = 'womens football pads'
$ws.Cells.Item(97, 1).Value = 'stretch mark men'
$ws.Cells.Item(98, 1).Value = 'youth xl baseball pants'
$ws.Cells.Item(99, 1).Value = 'youth basketball compression sleeve'
$ws.Cells.Item(100, 1).Value = 'sliding knee sleeve'
